# Applies the latest cryptocurrency price/volume snapshot to Sheet1.
# Column D = Price, Column E = Volume(1h) change, for each coin row (2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '20.103.43'
$ws.Range('E2').Value = '  -1.63%  '
# Row 3: Ethereum
$ws.Range('D3').Value = '1.422.55'
$ws.Range('E3').Value = '  -1.49%  '
# Row 4: TetherUSD
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  -0.28%  '
# Row 5: USDC
$ws.Range('E5').Value = '  -0.61%  '
# Row 6: BNB
$ws.Range('D6').Value = '''277.01'
$ws.Range('E6').Value = '  -0.37%  '
# Row 7: XRP
$ws.Range('E7').Value = '  -1.00%  '
# Row 8: Cardano
$ws.Range('D8').Value = '''0.3147'
$ws.Range('E8').Value = '  +2.35%  '
# Row 9: OKB
$ws.Range('D9').Value = '''39.49'
$ws.Range('E9').Value = '  -2.61%  '
# Row 10: Polygon
$ws.Range('D10').Value = '''1.061'
$ws.Range('E10').Value = '  +4.59%  '
# Row 11: Dogecoin
$ws.Range('D11').Value = '''0.06555'
$ws.Range('E11').Value = '  -0.21%  '
# Row 12: BinanceUSD
$ws.Range('D12').Value = '''0.9975'
$ws.Range('E12').Value = '  -0.60%  '
# Row 13: Polkadot
$ws.Range('D13').Value = '''5.532'
$ws.Range('E13').Value = '  +2.88%  '
# Row 14: Solana
$ws.Range('D14').Value = '''17.91'
$ws.Range('E14').Value = '  +3.70%  '
# Row 15: Chainlink
$ws.Range('D15').Value = '''6.204'
$ws.Range('E15').Value = '  +1.13%  '
# Row 16: WrappedEther
$ws.Range('D16').Value = '1.422.97'
$ws.Range('E16').Value = '  -1.47%  '
# Row 17: ShibaInu
$ws.Range('D17').Value = '''0.00001023'
$ws.Range('E17').Value = '  +1.25%  '
# Row 18: TRON
$ws.Range('D18').Value = '''0.05716'
$ws.Range('E18').Value = '  -2.61%  '
# Row 19: Dai
$ws.Range('D19').Value = '''0.9974'
$ws.Range('E19').Value = '  -0.57%  '
# Row 20: Litecoin
$ws.Range('D20').Value = '''71.71'
$ws.Range('E20').Value = '  -5.93%  '
# Row 21: Uniswap
$ws.Range('D21').Value = '''5.618'
$ws.Range('E21').Value = '  -1.98%  '
# Row 22: Avalanche
$ws.Range('D22').Value = '''14.88'
$ws.Range('E22').Value = '  +3.22%  '
# Row 23: Cosmos
$ws.Range('D23').Value = '''11.08'
$ws.Range('E23').Value = '  +1.82%  '
# Row 24: Toncoin
$ws.Range('D24').Value = '''2.225'
$ws.Range('E24').Value = '  -3.87%  '
# Row 25: WrappedBTC
$ws.Range('D25').Value = '20.139.12'
$ws.Range('E25').Value = '  -1.45%  '
# Row 26: LidoDAOToken
$ws.Range('D26').Value = '''2.293'
$ws.Range('E26').Value = '  +2.92%  '
# Row 27: Monero
$ws.Range('D27').Value = '''134.57'
$ws.Range('E27').Value = '  -6.07%  '
# Row 28: EthereumClassic
$ws.Range('E28').Value = '  +1.65%  '
# Row 29: WrappedliquidstakedEther2.0
$ws.Range('D29').Value = '1.581.64'
$ws.Range('E29').Value = '  -1.57%  '
# Row 30: BitcoinCash
$ws.Range('D30').Value = '''111.11'
$ws.Range('E30').Value = '  +1.46%  '
# Row 31: HuobiToken
$ws.Range('D31').Value = '''3.964'
$ws.Range('E31').Value = '  +4.98%  '
# Row 32: Filecoin
$ws.Range('D32').Value = '''5.289'
$ws.Range('E32').Value = '  -2.54%  '
# Row 33: ImmutableX
$ws.Range('D33').Value = '''0.8300'
$ws.Range('E33').Value = '  -8.40%  '
# Row 34: Stellar
$ws.Range('D34').Value = '''0.07823'
$ws.Range('E34').Value = '  +1.02%  '
# Row 35: WEMIXTOKEN
$ws.Range('D35').Value = '''1.478'
$ws.Range('E35').Value = '  +7.90%  '
# Row 36: InternetComputer(DFINITY)
$ws.Range('D36').Value = '''4.927'
$ws.Range('E36').Value = '  +4.11%  '
# Row 37: Hedera
$ws.Range('D37').Value = '''0.05872'
$ws.Range('E37').Value = '  +4.06%  '
# Row 38: FraxShare
$ws.Range('D38').Value = '''8.013'
$ws.Range('E38').Value = '  -3.52%  '
# Row 39: Frax
$ws.Range('D39').Value = '''0.9968'
$ws.Range('E39').Value = '  -0.58%  '
# Row 40: Aptos
$ws.Range('E40').Value = '  -2.46%  '
# Row 41: VeChain
$ws.Range('D41').Value = '''0.02067'
$ws.Range('E41').Value = '  +1.00%  '
# Row 42: TrustWalletToken
$ws.Range('D42').Value = '''1.111'
$ws.Range('E42').Value = '  -2.89%  '
# Row 43: Algorand
$ws.Range('D43').Value = '''0.1876'
$ws.Range('E43').Value = '  -2.13%  '
# Row 44: TheSandbox
$ws.Range('D44').Value = '''0.5349'
$ws.Range('E44').Value = '  +0.34%  '
# Row 45: EnergySwap
$ws.Range('E45').Value = '  +1.80%  '
# Row 46: PancakeSwap
$ws.Range('D46').Value = '''3.552'
$ws.Range('E46').Value = '  -1.16%  '
# Row 47: Quant
$ws.Range('D47').Value = '''118.48'
$ws.Range('E47').Value = '  +6.17%  '
# Row 48: Decentraland
$ws.Range('D48').Value = '''0.5239'
$ws.Range('E48').Value = '  +1.53%  '
# Row 49: NEARProtocol
$ws.Range('E49').Value = '  -0.41%  '
# Row 50: EOS
$ws.Range('E50').Value = '  -1.34%  '
# Row 51: PaxDollar
$ws.Range('D51').Value = '''0.9970'
$ws.Range('E51').Value = '  -0.63%  '
